$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits after the
#    "... email zu namen(teilnehmer aus club)" paragraph.
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) Add a new bold paragraph right after the
#    "testen gesamtes System (...) DIENSTAG" paragraph, and move the
#    _GoBack bookmark to the end of that new paragraph (login component
#    frontend task note).
# ------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("testen gesamtes System (Richtigkeit der Daten) bis spätestens DIENSTAG") | Out-Null
$targetPara = $findRng.Paragraphs(1)

$insertRng = $targetPara.Range
$insertRng.Collapse(0)
$insertRng.InsertParagraphAfter()

$newPara = $targetPara.Next()
$newRng = $newPara.Range

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>all-clubs-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>reduced</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>dataservice</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> implementieren (alle </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>clubs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> mitschicken die keinen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>verwalter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> haben=</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newRng.InsertXML($xmlFrag)
